# Update NATMI LR-pair output (Fgf2-Fgfr2) with newly re-computed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.456404
$ws.Range("H2").Value = 1.369212
$ws.Range("I2").Value = 0.01914960767004715
$ws.Range("J2").Value = 0.01914960767004715
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.036942
$ws.Range("N2").Value = 0.110826
$ws.Range("O2").Value = 0.02099032928903418
$ws.Range("P2").Value = 0.02099032928903418
$ws.Range("Q2").Value = 0.016860476568
$ws.Range("R2").Value = 0.151744289112
$ws.Range("S2").Value = 0.0004019565707501043
$ws.Range("T2").Value = 0.0004019565707501043

# Row 3 (ECs -> FAPs)
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.456404
$ws.Range("H3").Value = 1.369212
$ws.Range("I3").Value = 0.01914960767004715
$ws.Range("J3").Value = 0.01914960767004715
$ws.Range("O3").Value = 0.5358731102718634
$ws.Range("P3").Value = 0.5358731102718634
$ws.Range("Q3").Value = 0.4304399371133334
$ws.Range("R3").Value = 3.87395943402
$ws.Range("S3").Value = 0.0102617598226341
$ws.Range("T3").Value = 0.0102617598226341

# Row 4 (ECs -> MuSCs)
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.456404
$ws.Range("H4").Value = 1.369212
$ws.Range("I4").Value = 0.01914960767004715
$ws.Range("J4").Value = 0.01914960767004715
$ws.Range("O4").Value = 0.4431365604391025
$ws.Range("P4").Value = 0.4431365604391026
$ws.Range("Q4").Value = 0.3559493274653334
$ws.Range("R4").Value = 3.203543947188
$ws.Range("S4").Value = 0.008485891276662952
$ws.Range("T4").Value = 0.008485891276662952

# Row 5 (FAPs -> ECs)
$ws.Range("I5").Value = 0.8285024587002443
$ws.Range("J5").Value = 0.8285024587002443
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.036942
$ws.Range("N5").Value = 0.110826
$ws.Range("O5").Value = 0.02099032928903418
$ws.Range("P5").Value = 0.02099032928903418
$ws.Range("Q5").Value = 0.729463837178
$ws.Range("R5").Value = 6.565174534602001
$ws.Range("S5").Value = 0.01739053942489257
$ws.Range("T5").Value = 0.01739053942489257

# Row 6 (FAPs -> FAPs)
$ws.Range("I6").Value = 0.8285024587002443
$ws.Range("J6").Value = 0.8285024587002443
$ws.Range("O6").Value = 0.5358731102718634
$ws.Range("P6").Value = 0.5358731102718634
$ws.Range("S6").Value = 0.4439721894115859
$ws.Range("T6").Value = 0.4439721894115859

# Row 7 (FAPs -> MuSCs)
$ws.Range("I7").Value = 0.8285024587002443
$ws.Range("J7").Value = 0.8285024587002443
$ws.Range("O7").Value = 0.4431365604391025
$ws.Range("P7").Value = 0.4431365604391026
$ws.Range("S7").Value = 0.3671397298637659
$ws.Range("T7").Value = 0.3671397298637659

# Row 8 (MuSCs -> ECs)
$ws.Range("I8").Value = 0.1523479336297086
$ws.Range("J8").Value = 0.1523479336297086
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.036942
$ws.Range("N8").Value = 0.110826
$ws.Range("O8").Value = 0.02099032928903418
$ws.Range("P8").Value = 0.02099032928903418
$ws.Range("Q8").Value = 0.134136365058
$ws.Range("R8").Value = 1.207227285522
$ws.Range("S8").Value = 0.003197833293391508
$ws.Range("T8").Value = 0.003197833293391509

# Row 9 (MuSCs -> FAPs)
$ws.Range("I9").Value = 0.1523479336297086
$ws.Range("J9").Value = 0.1523479336297086
$ws.Range("O9").Value = 0.5358731102718634
$ws.Range("P9").Value = 0.5358731102718634
$ws.Range("S9").Value = 0.08163916103764338
$ws.Range("T9").Value = 0.08163916103764338

# Row 10 (MuSCs -> MuSCs)
$ws.Range("I10").Value = 0.1523479336297086
$ws.Range("J10").Value = 0.1523479336297086
$ws.Range("O10").Value = 0.4431365604391025
$ws.Range("P10").Value = 0.4431365604391026
$ws.Range("S10").Value = 0.06751093929867376
$ws.Range("T10").Value = 0.06751093929867377
